$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing values (second log entry for 03/16/2020 @ 19:30) ---
$ws.Range("B3").Value = 84
$ws.Range("B4").Value = 1000

# Row 5's "CCR2" label (D5) no longer applies to this entry - clear it.
$ws.Range("D5").ClearContents()

# Row 7 (Pulse) now also documents the CCR2 register, like row 4/5 do for ARR.
$ws.Range("D7").Value = "CCR2"
$ws.Range("D7").Font.Bold = $true

# --- New lookup table starting at row 10: Desired Frequency (Hz) / ARR / CCR2 ---
$ws.Range("A10").Value = "Desired Frequency (Hz)"
$ws.Range("B10").Value = "ARR"
$ws.Range("C10").Value = "CCR2"

$headerRange = $ws.Range("A10:C10")
$headerRange.Font.Bold = $true
$headerRange.Borders.Color = 0
$headerRange.Borders.LineStyle = 1

$freqValues = @(1,5,10,20,30,40,50,60,70,80,90,100,200,300,400,500,600,700,800,900,1000,1100,1200,1300,1400,1500,1600,1700,1800,1900,2000,2500,3000,3500,4000,4500,5000,6000,7000,8000,9000,10000)
for ($i = 0; $i -lt $freqValues.Length; $i++) {
    $row = 11 + $i
    $ws.Cells.Item($row, 1).Value = $freqValues[$i]
}

$dataRange = $ws.Range("A11:C52")
$dataRange.Borders.Color = 0
$dataRange.Borders.LineStyle = 1

# --- Cosmetic tweaks ---
$ws.Columns.Item(1).ColumnWidth = 21.14
$ws.Range("B11").Select()
